# Natmi following Dr Hou advice
#
# The sending-cluster/target-cluster grid for the Ntn1-Mcam ligand-receptor
# pair now includes an additional "ECs" cluster (alongside the existing
# FAPs, M2 and sCs clusters). This turns the previous 3 (sending) x 4
# (target) = 12-row grid into a full 4x4 = 16-row grid (rows 2-17).
# Ligand symbol is always "Ntn1" and Receptor symbol is always "Mcam".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{row=2; A="ECs"; B="Ntn1"; C="Mcam"; D="ECs"; E=2; F=0.6666666666666666; G=0.8750386666666667; H=2.625116; I=0.05304058862308838; J=0.05304058862308838; K=3; L=1; M=40.23443200000001; N=120.703296; O=0.5194057602668869; P=0.5194057602668869; Q=35.20668373137067; R=316.8601535823361; S=0.02754958725877841; T=0.02754958725877841},
    @{row=3; A="ECs"; B="Ntn1"; C="Mcam"; D="FAPs"; E=2; F=0.6666666666666666; G=0.8750386666666667; H=2.625116; I=0.05304058862308838; J=0.05304058862308838; K=3; L=1; M=2.585148666666667; N=7.755446; O=0.03337293561427507; P=0.03337293561427508; Q=2.262105042415111; R=20.358945381736; S=0.00177012014906158; T=0.00177012014906158},
    @{row=4; A="ECs"; B="Ntn1"; C="Mcam"; D="M2"; E=2; F=0.6666666666666666; G=0.8750386666666667; H=2.625116; I=0.05304058862308838; J=0.05304058862308838; K=3; L=1; M=0.525608; N=1.576824; O=0.006785328120013172; P=0.006785328120013173; Q=0.4599273235093334; R=4.139345911584001; S=0.0003598977974862923; T=0.0003598977974862924},
    @{row=5; A="ECs"; B="Ntn1"; C="Mcam"; D="sCs"; E=2; F=0.6666666666666666; G=0.8750386666666667; H=2.625116; I=0.05304058862308838; J=0.05304058862308838; K=3; L=1; M=34.11724066666667; N=102.351722; O=0.4404359759988248; P=0.4404359759988249; Q=29.85390478330578; R=268.6851430497521; S=0.02336098341776209; T=0.0233609834177621},
    @{row=6; A="FAPs"; B="Ntn1"; C="Mcam"; D="ECs"; E=3; F=1; G=10.61942; H=31.85826; I=0.6436975977089742; J=0.6436975977089742; K=3; L=1; M=40.23443200000001; N=120.703296; O=0.5194057602668869; P=0.5194057602668869; Q=427.26633186944; R=3845.396986824961; S=0.3343402401199985; T=0.3343402401199985},
    @{row=7; A="FAPs"; B="Ntn1"; C="Mcam"; D="FAPs"; E=3; F=1; G=10.61942; H=31.85826; I=0.6436975977089742; J=0.6436975977089742; K=3; L=1; M=2.585148666666667; N=7.755446; O=0.03337293561427507; P=0.03337293561427508; Q=27.45277945377333; R=247.07501508396; S=0.02148207848340513; T=0.02148207848340514},
    @{row=8; A="FAPs"; B="Ntn1"; C="Mcam"; D="M2"; E=3; F=1; G=10.61942; H=31.85826; I=0.6436975977089742; J=0.6436975977089742; K=3; L=1; M=0.525608; N=1.576824; O=0.006785328120013172; P=0.006785328120013173; Q=5.581652107359999; R=50.23486896624; S=0.004367699410519629; T=0.00436769941051963},
    @{row=9; A="FAPs"; B="Ntn1"; C="Mcam"; D="sCs"; E=3; F=1; G=10.61942; H=31.85826; I=0.6436975977089742; J=0.6436975977089742; K=3; L=1; M=34.11724066666667; N=102.351722; O=0.4404359759988248; P=0.4404359759988249; Q=362.3053078804134; R=3260.747770923721; S=0.2835075796950509; T=0.283507579695051},
    @{row=10; A="M2"; B="Ntn1"; C="Mcam"; D="ECs"; E=3; F=1; G=0.4029073333333333; H=1.208722; I=0.02442228319117198; J=0.02442228319117198; K=3; L=1; M=40.23443200000001; N=120.703296; O=0.5194057602668869; P=0.5194057602668869; Q=16.21074770530133; R=145.896729347712; S=0.01268507456836389; T=0.01268507456836389},
    @{row=11; A="M2"; B="Ntn1"; C="Mcam"; D="FAPs"; E=3; F=1; G=0.4029073333333333; H=1.208722; I=0.02442228319117198; J=0.02442228319117198; K=3; L=1; M=2.585148666666667; N=7.755446; O=0.03337293561427507; P=0.03337293561427508; Q=1.041575355556889; R=9.374178200011999; S=0.0008150432844925748; T=0.000815043284492575},
    @{row=12; A="M2"; B="Ntn1"; C="Mcam"; D="M2"; E=3; F=1; G=0.4029073333333333; H=1.208722; I=0.02442228319117198; J=0.02442228319117198; K=3; L=1; M=0.525608; N=1.576824; O=0.006785328120013172; P=0.006785328120013173; Q=0.2117713176586666; R=1.905941858928; S=0.0001657132048919842; T=0.0001657132048919843},
    @{row=13; A="M2"; B="Ntn1"; C="Mcam"; D="sCs"; E=3; F=1; G=0.4029073333333333; H=1.208722; I=0.02442228319117198; J=0.02442228319117198; K=3; L=1; M=34.11724066666667; N=102.351722; O=0.4404359759988248; P=0.4404359759988249; Q=13.74608645769822; R=123.714778119284; S=0.01075645213342352; T=0.01075645213342353},
    @{row=14; A="sCs"; B="Ntn1"; C="Mcam"; D="ECs"; E=3; F=1; G=4.600163333333334; H=13.80049; I=0.2788395304767656; J=0.2788395304767655; K=3; L=1; M=40.23443200000001; N=120.703296; O=0.5194057602668869; P=0.5194057602668869; Q=185.0849588238934; R=1665.76462941504; S=0.1448308583197462; T=0.1448308583197462},
    @{row=15; A="sCs"; B="Ntn1"; C="Mcam"; D="FAPs"; E=3; F=1; G=4.600163333333334; H=13.80049; I=0.2788395304767656; J=0.2788395304767655; K=3; L=1; M=2.585148666666667; N=7.755446; O=0.03337293561427507; P=0.03337293561427508; Q=11.89210610761556; R=107.02895496854; S=0.00930569369731579; T=0.00930569369731579},
    @{row=16; A="sCs"; B="Ntn1"; C="Mcam"; D="M2"; E=3; F=1; G=4.600163333333334; H=13.80049; I=0.2788395304767656; J=0.2788395304767655; K=3; L=1; M=0.525608; N=1.576824; O=0.006785328120013172; P=0.006785328120013173; Q=2.417882649306667; R=21.76094384376; S=0.001892017707115267; T=0.001892017707115267},
    @{row=17; A="sCs"; B="Ntn1"; C="Mcam"; D="sCs"; E=3; F=1; G=4.600163333333334; H=13.80049; I=0.2788395304767656; J=0.2788395304767655; K=3; L=1; M=34.11724066666667; N=102.351722; O=0.4404359759988248; P=0.4404359759988249; Q=156.9448795493089; R=1412.50391594378; S=0.1228109607525883; T=0.1228109607525883}
)

$cols = "A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T"

foreach ($r in $data) {
    $rowNum = $r.row
    foreach ($col in $cols) {
        $ws.Range("$col$rowNum").Value = $r[$col]
    }
}

Write-Host "Wrote $($data.Count) data rows (rows 2..$(1 + $data.Count))"
